$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like values (e.g. "24.681.17", "0.07990") are stored as text,
# matching the original inlineStr cell type, rather than being auto-converted
# to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.681.17"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.694.93"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.31"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3961"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4077"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.496"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.10"
$ws.Range("E11").Value = "  -5.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08925"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.202"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.62"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.189"
$ws.Range("E15").Value = "  +10.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001327"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.693.08"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.25"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07010"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.74"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.035"
$ws.Range("E21").Value = "  +6.77%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.30"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.672.80"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.208"
$ws.Range("E25").Value = "  +5.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.346"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.69"
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.24"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "137.10"
$ws.Range("E29").Value = "  +4.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.174"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.509"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.882.51"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08621"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.165"
$ws.Range("E35").Value = "  -6.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.50"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2744"
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.922"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.49"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09191"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02732"
$ws.Range("E41").Value = "  +7.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.478"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7681"
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.09"
$ws.Range("E44").Value = "  +6.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.627"
$ws.Range("E45").Value = "  +8.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7195"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.231"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.63"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.327"
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07990"
$ws.Range("E51").Value = "  +1.80%  "

Write-Host "Updated 97 cells"
